# The edit permutes the data rows 2-10 (entire row contents, columns A:AY)
# according to the mapping below: sourceRow -> destinationRow.
# (Row content that used to live at source row now lives at destination row.)
#
#   before row  9 -> after row 2
#   before row  2 -> after row 3
#   before row  3 -> after row 4
#   before row  4 -> after row 5
#   before row  5 -> after row 6
#   before row  6 -> after row 7
#   before row 10 -> after row 8
#   before row  7 -> after row 9
#   before row  8 -> after row 10
#
# i.e. two row cycles: (9 2 3 4 5 6 7) and (10 8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$moveMap = @{
    9  = 2
    2  = 3
    3  = 4
    4  = 5
    5  = 6
    6  = 7
    10 = 8
    7  = 9
    8  = 10
}

$firstCol = "A"
$lastCol  = "AY"
$tempBase = 200   # scratch rows, far below the real data, used for staging

# Step 1: snapshot every affected row (as a whole, A:AY) into a temporary
# holding row using Copy (preserves types/formatting exactly - e.g. text
# cells that look like dates stay text instead of being re-parsed).
# We must stage everything first since the row move is a permutation
# (made of two cycles), so a direct row-to-row copy would clobber a row
# before it has been read.
$tempRowOf = @{}
$i = 0
foreach ($srcRow in $moveMap.Keys) {
    $tempRow = $tempBase + $i
    $tempRowOf[$srcRow] = $tempRow

    $tempRange = $ws.Range("$firstCol$tempRow`:$lastCol$tempRow")
    $tempRange.Clear()
    $ws.Range("$firstCol$srcRow`:$lastCol$srcRow").Copy($tempRange)

    $i++
}

# Step 2: copy each staged row into its real destination row.
foreach ($srcRow in $moveMap.Keys) {
    $dstRow  = $moveMap[$srcRow]
    $tempRow = $tempRowOf[$srcRow]

    $dstRange = $ws.Range("$firstCol$dstRow`:$lastCol$dstRow")
    $dstRange.Clear()
    $ws.Range("$firstCol$tempRow`:$lastCol$tempRow").Copy($dstRange)
}

# Step 3: wipe the scratch rows so the sheet dimensions/content stay clean.
foreach ($srcRow in $moveMap.Keys) {
    $tempRow = $tempRowOf[$srcRow]
    $ws.Range("$firstCol$tempRow`:$lastCol$tempRow").Clear()
}
